# The deck ships with two slides:
#   1. SlideID 256 - cover/title slide            -> kept
#   2. SlideID 257 - "Danh sach thanh vien:" (member-list placeholder slide) -> removed
#
# The canonical diff drops the second <p:sldId> entry from the
# presentation's slide list and removes ppt/slides/slide2.xml (plus its
# relationship) entirely. Deleting the Slide object through the
# PowerPoint object model reproduces exactly that cascade: the slide
# part, its rels, the [Content_Types].xml override and the
# <p:sldIdLst> entry are all cleaned up together by PowerPoint itself.

$p = $ppt.ActivePresentation

$slideToRemove = $p.Slides.Item($p.Slides.Count)
$slideToRemove.Delete()
